$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new row of data (row 6) — a newly solved problem entry,
# mirroring the existing table layout/categories.
$ws.Range("A6").Value = "D&C"
$ws.Range("B6").Value = "Convert Sorted Array to Binary Search Tree"
$ws.Range("C6").Value = "d&c, recursion;"

# Match the cell formatting used by the other "category" rows (B2:B4)
# by copying their format onto the new cell.
$ws.Range("B2").Copy()
$ws.Range("B6").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Update the selection to match the authored state.
$ws.Range("C6").Select()
